# Weekly update: prepend a new week's worth of data (2 rows) for
# Hortaliza, Terminal La Palmera de La Serena - Brocoli.
# This pushes the existing data block (rows 509:542) down by two rows
# (to 511:544) and fills the freed rows 509:510 with the new week's
# figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 509 downward by inserting two fresh rows above the block.
$ws.Rows("509:510").Insert()

# --- Row 509 (Primera) ---
$ws.Range("A509").Value = 8
$ws.Range("B509").Value = "Terminal La Palmera de La Serena"
$ws.Range("C509").Value = "Coquimbo"
$ws.Range("D509").Value = 44585
$ws.Range("E509").Value = 4
$ws.Range("F509").Value = 100112023
$ws.Range("G509").Value = "Brócoli"
$ws.Range("H509").Value = "Sin especificar"
$ws.Range("I509").Value = "Primera"
$ws.Range("J509").Value = 2000
$ws.Range("K509").Value = 750
$ws.Range("L509").Value = 800
$ws.Range("M509").Value = 775
$ws.Range("N509").Value = "$/unidad"
$ws.Range("O509").Value = "Provincia del Elquí"
$ws.Range("P509").Value = 775
$ws.Range("Q509").Value = 1
$ws.Range("R509").Value = "Hortaliza"

# --- Row 510 (Segunda) ---
$ws.Range("A510").Value = 8
$ws.Range("B510").Value = "Terminal La Palmera de La Serena"
$ws.Range("C510").Value = "Coquimbo"
$ws.Range("D510").Value = 44585
$ws.Range("E510").Value = 4
$ws.Range("F510").Value = 100112023
$ws.Range("G510").Value = "Brócoli"
$ws.Range("H510").Value = "Sin especificar"
$ws.Range("I510").Value = "Segunda"
$ws.Range("J510").Value = 1600
$ws.Range("K510").Value = 650
$ws.Range("L510").Value = 700
$ws.Range("M510").Value = 675
$ws.Range("N510").Value = "$/unidad"
$ws.Range("O510").Value = "Provincia del Elquí"
$ws.Range("P510").Value = 675
$ws.Range("Q510").Value = 1
$ws.Range("R510").Value = "Hortaliza"
